$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Email" header cell's dependent value: C2 currently has a hyperlink
# display text "your email" shared string that must become "sender mail".
$ws.Range("C2").Value = "sender mail"

# Remove the hyperlink on C2 (keep the text, drop the mailto: link).
$ws.Hyperlinks.Delete()

# Add a new empty styled cell at C3, matching the Hyperlink style used by C2,
# and select it.
$ws.Range("C3").Value = ""
$ws.Range("C3").Style = $ws.Range("C2").Style

$ws.Range("C2").Select()
